# Update DateBase/orders/Fresh bloom Flowers_2025-10-7.xlsx
# - Orders sheet: extend order table with a new group (#13 and #14), grows
#   from row 51 to row 61; F51 quantity corrected from 2 -> 20.
# - Summary sheet: G2 running-tally string gets a new chunk appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# The PackageID (A) and Number (F) columns hold digit strings stored as TEXT
# (not numbers) throughout this sheet. Pre-format the cells we are about to
# touch as Text so Excel's auto-detection doesn't coerce these into numeric
# values (which would also strip the existing "store as text" look of the
# rest of the column). Multi-area (comma) ranges aren't reliable here, so
# each contiguous block is formatted separately.
$ws.Range("F51:F60").NumberFormat = "@"
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A58").NumberFormat = "@"

# Fix the quantity on the existing last row (577_腊梅白 / wax white), 2 -> 20
$ws.Cells.Item(51, 6).Value = "20"

# New data rows appended after row 51 (columns: A=PackageID, C=FlowerName, F=Number)
$newRows = @(
    @{ Row = 52; A = $null; C = "579_腊梅红_wax red_undefined_1bunch"; F = "5" },
    @{ Row = 53; A = "13";  C = "2_粉洋桔梗_Pink Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "15" },
    @{ Row = 54; A = $null; C = "1_白洋桔梗_White Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "30" },
    @{ Row = 55; A = $null; C = "12_肉粉洋桔梗_Peach Lisianthus_Eustoma grandiflorum (Raf.) Shinners_800/600g"; F = "15" },
    @{ Row = 56; A = $null; C = "300_白星_White Gypso_ gypsophila_1kg"; F = "9" },
    @{ Row = 57; A = $null; C = "686_百合小粉仙_undefined_undefined_1bunch"; F = "5" },
    @{ Row = 58; A = "14";  C = "106_绣球单瓣粉_Hydrangea Pink S_Hydrangea L._1stem"; F = "20" },
    @{ Row = 59; A = $null; C = "798_朱玉大菊_undefined_undefined_5stems"; F = "5" },
    @{ Row = 60; A = $null; C = "535_雪果白_snow berry white_undefined_1bunch"; F = "5" },
    @{ Row = 61; A = $null; C = "780_贝壳草_undefined_undefined_1bunch"; F = $null }
)

foreach ($r in $newRows) {
    if ($null -ne $r.A) {
        $ws.Cells.Item($r.Row, 1).Value = $r.A
    }
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($null -ne $r.F) {
        $ws.Cells.Item($r.Row, 6).Value = $r.F
    }
}

# Update the Summary sheet running-tally string (G2) with the appended chunk.
# Force text format first so the long digit string isn't coerced into a
# double (which would lose precision / truncate the value).
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2, 7).NumberFormat = "@"
$ws2.Cells.Item(2, 7).Value = "013242011.5111253551013822121431751240503325531555211091029323555555554512520102051530159520550"
